# Fix for multiple result
# The shared-string entry "4709 Gonzales St #A" was a duplicate/incorrect
# address; correct it to "4709 Gonzales St". Also move the active
# selection from I21 to A21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the street address text in A19 (was "4709 Gonzales St #A")
$ws.Range("A19").Value = "4709 Gonzales St"

# Update the sheet's active selection/cell to A21
$ws.Range("A21").Select() | Out-Null
